$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): corrected B3 value; D3 no longer carries a value
$ws.Range("B3").Value = 6137285.521024222
$ws.Range("D3").ClearContents()

# Row 4 (Methanol): corrected C4 value
$ws.Range("C4").Value = 4692.345966963499

# Row 5 (Ammonia): corrected C5 value
$ws.Range("C5").Value = 12413.92129710552

# Row 7 used to be "Other" - it is now "Biogas" with a corrected value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 1069.278320203996

# A new row 8 holds the (renamed) "Other" category with its own value.
# Copy row 7's formatting down to row 8 first so the label cell picks up
# the same bold/border/center style used by the other row headers.
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)

$ws.Range("A8").Value = "Other"
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = 85.62302098320228
